# msz - field hint and error checks part 4 + dialog resolve
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert 3 new rows before row 7 so the single "Goto ..." 2-row block
# (rows 5-6) becomes a 5-row block (rows 5-9), shifting everything
# below (old rows 7-16) down to new rows 10-19.
$ws.Range("A7:H9").EntireRow.Insert()

# Clear the two old "Goto ... page" rows first (the Insert above only
# pushed down rows 7+, rows 5-6 keep their original cells otherwise).
$ws.Range("A5:H6").ClearContents()

# Rewrite the "Goto ... page" block (rows 5-9): one row per dialog page,
# with the <SELECT> marker stepping one column to the right each row.
# (Row labels are entered in this order so new shared-string entries are
# minted in the same sequence as the source workbook.)
$ws.Range("A8").Value = "Goto price option page"
$ws.Range("A5").Value = "Goto vehicle data page"
$ws.Range("A6").Value = "Goto insurant data page"
$ws.Range("A7").Value = "Goto product data page"
$ws.Range("A9").Value = "Goto send quote page"

$ws.Range("B5").Value = "<SET>"
$ws.Range("C5").Value = "<SELECT>"
$ws.Range("H5").Value = "<NOP>"

$ws.Range("B6").Value = "<SET>"
$ws.Range("D6").Value = "<SELECT>"
$ws.Range("H6").Value = "<NOP>"

$ws.Range("B7").Value = "<SET>"
$ws.Range("E7").Value = "<SELECT>"
$ws.Range("H7").Value = "<NOP>"

$ws.Range("B8").Value = "<SET>"
$ws.Range("F8").Value = "<SELECT>"
$ws.Range("H8").Value = "<NOP>"

$ws.Range("B9").Value = "<SET>"
$ws.Range("G9").Value = "<SELECT>"
$ws.Range("H9").Value = "<NOP>"

# Rows 10-13 keep the highlighted (yellow) style of the old rows 7-10,
# content unchanged, just shifted down by 3 rows by the insert above.
$ws.Range("A14").Value = "Button Next from Page VehicleData"
$ws.Range("C14").Value = "Button Next"

$ws.Range("A19").Value = "Send Quote - Button Main Page"
$ws.Range("G19").Value = "Button Main Page"

$ws.Range("E12").Select()
